# Automatische test-sync: 2025-06-17 22:29:11
# Append two new rows to the "Logs" sheet and update the "Dashboard" summary counts.

$wb = $excel.ActiveWorkbook

# --- Update "Logs" sheet: append rows 53 and 54 ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A53").Value = "Sollicitatie marketingfunctie"
$logs.Range("B53").Value = "mailmind.test@zohomail.eu"
$logs.Range("C53").Value = "Hierbij solliciteer ik voor de functie van marketeer. Zie bijlage voor CV."
$logs.Range("D53").Value = "Overig"
$logs.Range("F53").Value = "2025-06-17 22:29:05"
$logs.Range("G53").Value = "Nee"

$logs.Range("A54").Value = "Offerte voor zakelijke samenwerking"
$logs.Range("B54").Value = "mailmind.test@zohomail.eu"
$logs.Range("C54").Value = "Kunt u mij een offerte sturen voor 100 stuks product X?"
$logs.Range("D54").Value = "Bestelling"
$logs.Range("F54").Value = "2025-06-17 22:29:05"
$logs.Range("G54").Value = "Nee"

# Extend the conditional formatting ranges so the new rows are covered too
# (FormatConditions.ModifyAppliesToRange must be invoked per-rule, not on the collection)
$catConditions = $logs.Range("D2:D52").FormatConditions
for ($i = 1; $i -le $catConditions.Count; $i++) {
    $catConditions.Item($i).ModifyAppliesToRange($logs.Range("D2:D54"))
}

$answeredConditions = $logs.Range("G2:G52").FormatConditions
for ($i = 1; $i -le $answeredConditions.Count; $i++) {
    $answeredConditions.Item($i).ModifyAppliesToRange($logs.Range("G2:G54"))
}

# --- Update "Dashboard" sheet: refresh category counts ---
$dashboard = $wb.Worksheets.Item("Dashboard")

$dashboard.Range("B3").Value = 14
$dashboard.Range("B6").Value = 4
